$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Clear the candidate data (name, mobile, application no, category, panel)
# for rows 2-15, leaving the existing cell formatting/styles intact.
$ws.Range("B2:F15").ClearContents()

# Move the active selection to H2 (matches the saved view in the workbook).
$ws.Range("H2").Select()
